# Update LDLC prices history
#
# A new price-snapshot column is inserted right before the existing
# "nom" / "url_produit" columns (BT/BU), pushing them one column to the
# right (-> BU/BV). The freshly inserted column is filled with a copy of
# the previous (rightmost) price column (BS), mirroring how every other
# timestamped price column in this sheet was produced, and its header
# (row 1) gets the new snapshot's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at BT; old BT (nom) -> BU, old BU (url_produit) -> BV.
$ws.Columns("BT").Insert()

# 2. Fill the new BT column (rows 2-206) with a copy of column BS (the
#    last existing price snapshot), so each row keeps the last known
#    price/blank exactly like BS already has it.
$ws.Range("BS2:BS206").Copy()
$ws.Range("BT2:BT206").PasteSpecial()

# 3. Set the new column's header to the new snapshot timestamp.
$ws.Range("BT1").Value = "2026-01-30 22:12:56"
